# Writing value function: allows writing a value or a formula/instruction
# into a given cell position, on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Write-CellValue {
    param($Sheet, [string]$Cell, $Value)

    $range = $Sheet.Range($Cell)
    if ($null -ne $Value -and $Value -is [string] -and $Value.StartsWith("=")) {
        $range.Formula = $Value
    }
    else {
        $range.Value = $Value
    }
}

# --- Update the number formatting of the existing "note" column cells so
#     every mark is displayed with two decimal places (0.00) ---
$noteCells = "B4", "B7", "B8", "B10", "B12", "B17", "B18"
foreach ($cellRef in $noteCells) {
    $ws.Range($cellRef).NumberFormat = "0.00"
}

# --- New cell D9: styled placeholder (green Consolas font) ready to host
#     a future instruction, no value written yet ---
$d9 = $ws.Range("D9")
$d9.Font.Name = "Consolas"
$d9.Font.Size = 11
$d9.Font.Color = 7979928
$d9.NumberFormat = "0.00"

# --- Write a plain value in a given cell position ---
Write-CellValue $ws "H5" 8

# --- Write an instruction (formula) in a given cell position, on the file
#     whose path was given too (the already opened workbook) ---
Write-CellValue $ws "F6" "=SUM(B3:B20)"
Write-CellValue $ws "H6" "=SUM(B3:B20)"

# --- Selection / active cell ---
[void]$ws.Range("H6").Select()

# --- Page setup ---
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
